$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 52, shifting existing rows 52:96 down to 53:97
$ws.Rows.Item(52).Insert()

# Populate the new row 52 with its data
$ws.Cells.Item(52, 1).Value = 11
$ws.Cells.Item(52, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(52, 3).Value = "Bíobío"
$ws.Cells.Item(52, 4).Value = 44484
$ws.Cells.Item(52, 5).Value = 8
$ws.Cells.Item(52, 6).Value = 100112003
$ws.Cells.Item(52, 7).Value = "Ajo"
$ws.Cells.Item(52, 8).Value = "Chino"
$ws.Cells.Item(52, 9).Value = "Primera"
$ws.Cells.Item(52, 10).Value = 430
$ws.Cells.Item(52, 11).Value = 14000
$ws.Cells.Item(52, 12).Value = 15000
$ws.Cells.Item(52, 13).Value = 14535
$ws.Cells.Item(52, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(52, 15).Value = "China"
$ws.Cells.Item(52, 16).Value = 1454
$ws.Cells.Item(52, 17).Value = 10
$ws.Cells.Item(52, 18).Value = "Hortaliza"
